$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Measures  (sheet1) -------------------------------------------------
# ---------------------------------------------------------------------------
$wsMeasures = $wb.Worksheets.Item("Measures")

# Remove the long tail of placeholder / leftover rows (26-187) that only
# contained stray formatting, leaving the real data in rows 1-25.
$wsMeasures.Range("A26:A187").EntireRow.Delete()

# Insert a brand new column before the existing "comment" column (O) so the
# old column O becomes P and a fresh column O is available for "item_num".
$wsMeasures.Columns.Item(15).Insert()
$wsMeasures.Range("O1").Value = "item_num"
$wsMeasures.Range("O2:O25").Value = 1

# Drop the now-unused explicit cell styles so the sheet reverts to the
# workbook default style.
$wsMeasures.Rows.Item(1).ClearFormats()
$wsMeasures.Range("I2:I25").ClearFormats()

$wsMeasures.Range("E6").Select()

# ---------------------------------------------------------------------------
# Sheet: ID  (sheet2) --------------------------------------------------------
# ---------------------------------------------------------------------------
$wsID = $wb.Worksheets.Item("ID")
$wsID.Columns.Item(15).Insert()
$wsID.Range("O1").Value = "item_num"
$wsID.Rows.Item(1).ClearFormats()
$wsID.Range("O1:O1048576").Select()

# ---------------------------------------------------------------------------
# Sheet: Dems  (sheet3) ------------------------------------------------------
# ---------------------------------------------------------------------------
$wsDems = $wb.Worksheets.Item("Dems")
$wsDems.Columns.Item(15).Insert()
$wsDems.Range("O1").Value = "item_num"
$wsDems.Rows.Item(1).ClearFormats()
$wsDems.Range("O1:O1048576").Select()

# ---------------------------------------------------------------------------
# Sheet: Dates  (sheet4) ------------------------------------------------------
# ---------------------------------------------------------------------------
$wsDates = $wb.Worksheets.Item("Dates")
$wsDates.Columns.Item(15).Insert()
$wsDates.Range("O1").Value = "item_num"
$wsDates.Rows.Item(1).ClearFormats()
$wsDates.Range("O1:O1048576").Select()

# ---------------------------------------------------------------------------
# Sheet: NewVars  (sheet5) ----------------------------------------------------
# ---------------------------------------------------------------------------
$wsNewVars = $wb.Worksheets.Item("NewVars")
$wsNewVars.Columns.Item(15).Insert()
$wsNewVars.Range("O1").Value = "item_num"
$wsNewVars.Range("O2:O5").Value = 5
$wsNewVars.Rows.Item(1).ClearFormats()

# NewVars becomes the active sheet/tab, with O6 as the last selection made -
# doing this last makes NewVars the workbook's active tab on save.
$wsNewVars.Range("O6").Select()

# ---------------------------------------------------------------------------
# Workbook level: update the hidden filter-database range for Measures so it
# covers the newly inserted column (O -> P).
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Measures!_FilterDatabase") {
        $n.RefersTo = "=Measures!`$A`$1:`$P`$28"
    }
}
